$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks cleanly (will be recreated after the row insert)
$ws.Hyperlinks.Delete()

# Insert two new rows at the top of the data (rows 2 and 3), pushing all other rows down
$ws.Rows("2:3").Insert()

# Refresh timestamp for every data row (2..12) to the new scrape time
$ws.Range("A2:A12").Value = "2025-12-18 18:28:38"

# Row 2
$ws.Range("B2").Value = "【急募】AI医療系請求IOSアプリ開発のフリーランス募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5456942"
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = "🔥AI,Ai ◆開発 ◇アプリ"

# Row 3
$ws.Range("B3").Value = "EC×AIプロダクト/業務改善リード"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5450024"
$ws.Range("G3").Value = 338
$ws.Range("H3").Value = "🔥AI,Ai ◇業務改善"

# Row 4
$ws.Range("B4").Value = "初回 既存システムのRuby、Ruby on Railsバージョンアップ及び追加改修"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5456434"
$ws.Range("G4").Value = 318
$ws.Range("H4").Value = "🔥AI,Ai"

# Row 5
$ws.Range("B5").Value = "【Zapier設定のみ!作業時間~2時間】スプレッドシート・Gドライブ自動化構築(設計済)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5456066"
$ws.Range("G5").Value = 255
$ws.Range("H5").Value = "🔥API ◆自動化"

# Row 6
$ws.Range("B6").Value = "【スマホアプリ開発】 音声データ推定アプリの依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5456360"
$ws.Range("G6").Value = 175
$ws.Range("H6").Value = "★スマホアプリ ◆開発 ◇アプリ"

# Row 7
$ws.Range("B7").Value = "【急募】多店舗パーソナルジム向け予約・顧客管理システム開発"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5456461"
$ws.Range("G7").Value = 160
$ws.Range("H7").Value = "◆開発,システム開発 ◇管理"

# Row 8
$ws.Range("B8").Value = "【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5456658"
$ws.Range("G8").Value = 158
$ws.Range("H8").Value = "◆自動化,スクレイピング ◇管理"

# Row 9
$ws.Range("B9").Value = "【完全在宅/時給1,400円】IT・業務効率化経験を活かせる!社内エンジニア兼総務スタッフを募集!"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5456452"
$ws.Range("G9").Value = 70
$ws.Range("H9").Value = "◆効率化"

# Row 10
$ws.Range("B10").Value = "【急募】データ活用インフラ要件整理のコンサルタント募集(1人月/月)"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5456545"
$ws.Range("G10").Value = 55
$ws.Range("H10").Value = "◆コンサル"

# Row 11
$ws.Range("B11").Value = "【急募】LINEシステム構築・保守運用のプロフェッショナルを求む!"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5456063"
$ws.Range("G11").Value = 33
$ws.Range("H11").ClearContents()

# Row 12
$ws.Range("B12").Value = "【準委任】音声データ収集プロジェクト/PM・ディレクター募集"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5456449"
$ws.Range("G12").Value = 25
$ws.Range("H12").ClearContents()

# Re-create hyperlinks for the URL column in row order, then restyle as Hyperlink
$ws.Hyperlinks.Add($ws.Range("F2"), $ws.Range("F2").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), $ws.Range("F3").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), $ws.Range("F4").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), $ws.Range("F5").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), $ws.Range("F6").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), $ws.Range("F7").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), $ws.Range("F8").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), $ws.Range("F9").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), $ws.Range("F10").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), $ws.Range("F11").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), $ws.Range("F12").Value2) | Out-Null
$ws.Range("F2:F12").Style = "Hyperlink"

$ws.Range("A1").Select()